# The document contains a single-column table. Each of the first 12 rows
# holds one statistic value in its own run, and rows 44-46 hold several
# tab-separated values packed into one run. Per the diff:
#  - Rows 1-3 become "0M" (memory stat now reported differently)
#  - Row 4's value changes from 701 to 2101
#  - Rows 5-11 get updated numeric values
#  - Rows 44-46 collapse from their long tab-separated payload down to a
#    single value (the one previously held by rows 1-3 before this edit)

$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "2101"
$t.Cell(5, 1).Range.Text  = "0.00001"
$t.Cell(6, 1).Range.Text  = "0.00053"
$t.Cell(7, 1).Range.Text  = "0.00016"
$t.Cell(8, 1).Range.Text  = "0.00002"
$t.Cell(9, 1).Range.Text  = "0.00032"
$t.Cell(10, 1).Range.Text = "0.00036"
$t.Cell(11, 1).Range.Text = "0.00038"
$t.Cell(12, 1).Range.Text = "0.42126"

$t.Cell(44, 1).Range.Text = "99.87"
$t.Cell(45, 1).Range.Text = "0.42"
$t.Cell(46, 1).Range.Text = "318"
